$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (pushes existing rows 9..120 down to 10..121,
# carrying the style of row 9 with it, same as a normal Excel row insert).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new price-record data.
# (Columns A, B, C, E, F, G, H, I, J, L, M, Q, R, T are identical to the
# record that used to occupy row 9 before the insert -- i.e. what is now
# row 10 -- so copy them across; only D, K, N, O, P, S are genuinely new.)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44616
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100109
$ws.Range("H9").Value = "Uva"
$ws.Range("I9").Value = 100109001
$ws.Range("J9").Value = "Uva"
$ws.Range("K9").Value = "Flame Seedless"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 10500
$ws.Range("Q9").Value = "$/bandeja 18 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 583
$ws.Range("T9").Value = 18
